$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'95.291.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -2.07%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.605.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.98%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'2.65"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +38.19%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  +0.05%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'222.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -6.02%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'634.56"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -3.69%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.414"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.94%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +9.37%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.07%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'3.604.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.94%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'47.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +5.57%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.210"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.08%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000288"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -9.55%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -7.22%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'4.279.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.93%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'95.250.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.65%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'8.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.16%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.66%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'3.607.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.02%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.541"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +6.06%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +46.65%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'509.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -3.11%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -7.69%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'119.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +13.15%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0000198"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -11.27%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'6.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -2.09%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.789.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'12.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -7.02%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'12.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.07%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.06%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.35%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.610"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.92%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'32.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -7.28%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.47%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'8.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.68%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'574.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -10.40%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'6.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'41.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.74%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.491"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.21%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0501"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +10.15%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.154"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -7.51%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.947"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.28%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.47%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'226.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +10.09%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.30%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'23.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.70%  "
$ws.Range("E51").Style = "Normal"
